$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D-column price values as text (Excel auto-converts plain numeric
# strings to Number on assignment, just like typing into a cell would; a
# leading apostrophe forces text entry, matching the original inlineStr type).

$ws.Range('D2').Value = "'24.896.95"
$ws.Range('E2').Value = '  +1.47%  '

$ws.Range('D3').Value = "'1.710.03"
$ws.Range('E3').Value = '  +1.49%  '

$ws.Range('D4').Value = "'1.003"
$ws.Range('E4').Value = '  -0.59%  '

$ws.Range('D5').Value = "'315.31"
$ws.Range('E5').Value = '  +0.45%  '

$ws.Range('E6').Value = '  -0.45%  '

$ws.Range('D7').Value = "'0.4040"
$ws.Range('E7').Value = '  +3.69%  '

$ws.Range('D8').Value = "'0.4057"
$ws.Range('E8').Value = '  +1.03%  '

$ws.Range('E9').Value = '  -0.56%  '

$ws.Range('E10').Value = '  -0.26%  '

$ws.Range('D11').Value = "'53.72"
$ws.Range('E11').Value = '  +1.50%  '

$ws.Range('D12').Value = "'0.08812"
$ws.Range('E12').Value = '  +1.50%  '

$ws.Range('D13').Value = "'26.24"
$ws.Range('E13').Value = '  +7.60%  '

$ws.Range('D14').Value = "'7.524"
$ws.Range('E14').Value = '  -1.43%  '

$ws.Range('D15').Value = "'8.005"
$ws.Range('E15').Value = '  +1.03%  '

$ws.Range('D16').Value = "'0.00001345"
$ws.Range('E16').Value = '  +1.23%  '

$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = "'95.59"
$ws.Range('E17').Value = '  -2.36%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = "'0.07176"
$ws.Range('E18').Value = '  +1.18%  '

$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = "'1.557.81"
$ws.Range('E19').Value = '  -7.51%  '

$ws.Range('D20').Value = "'21.05"
$ws.Range('E20').Value = '  +7.55%  '

$ws.Range('D21').Value = "'7.295"
$ws.Range('E21').Value = '  +0.70%  '

$ws.Range('D22').Value = "'1.004"
$ws.Range('E22').Value = '  -0.73%  '

$ws.Range('D23').Value = "'14.48"
$ws.Range('E23').Value = '  +2.52%  '

$ws.Range('D24').Value = "'24.878.87"
$ws.Range('E24').Value = '  +1.45%  '

$ws.Range('D25').Value = "'2.337"
$ws.Range('E25').Value = '  -0.56%  '

$ws.Range('D26').Value = "'2.889"
$ws.Range('E26').Value = '  -3.26%  '

$ws.Range('D27').Value = "'6.400"
$ws.Range('E27').Value = '  +22.18%  '

$ws.Range('D28').Value = "'23.09"
$ws.Range('E28').Value = '  +2.36%  '

$ws.Range('D29').Value = "'162.11"
$ws.Range('E29').Value = '  +0.60%  '

$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = "'144.10"
$ws.Range('E30').Value = '  +5.81%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = "'8.394"
$ws.Range('E31').Value = '  -0.74%  '

$ws.Range('D32').Value = "'2.282"
$ws.Range('E32').Value = '  +15.33%  '

$ws.Range('D33').Value = "'0.08755"
$ws.Range('E33').Value = '  +0.58%  '

$ws.Range('D34').Value = "'0.03184"
$ws.Range('E34').Value = '  +10.69%  '

$ws.Range('D35').Value = "'1.781.71"
$ws.Range('E35').Value = '  -4.82%  '

$ws.Range('D36').Value = "'7.206"
$ws.Range('E36').Value = '  -3.18%  '

$ws.Range('D37').Value = "'1.028"
$ws.Range('E37').Value = '  +0.11%  '

$ws.Range('D38').Value = "'0.2875"
$ws.Range('E38').Value = '  +6.58%  '

$ws.Range('D39').Value = "'0.8443"
$ws.Range('E39').Value = '  +9.77%  '

$ws.Range('D40').Value = "'10.84"
$ws.Range('E40').Value = '  +1.91%  '

$ws.Range('D41').Value = "'0.09472"
$ws.Range('E41').Value = '  +4.30%  '

$ws.Range('D42').Value = "'14.23"
$ws.Range('E42').Value = '  +1.77%  '

$ws.Range('D43').Value = "'1.479"
$ws.Range('E43').Value = '  +2.09%  '

$ws.Range('D44').Value = "'17.57"
$ws.Range('E44').Value = '  +6.50%  '

$ws.Range('D45').Value = "'2.725"
$ws.Range('E45').Value = '  +6.85%  '

$ws.Range('D46').Value = "'0.7444"
$ws.Range('E46').Value = '  +5.01%  '

$ws.Range('E47').Value = '  +0.95%  '

$ws.Range('D48').Value = "'1.378"
$ws.Range('E48').Value = '  +3.76%  '

$ws.Range('D49').Value = "'1.002"
$ws.Range('E49').Value = '  -0.77%  '

$ws.Range('D50').Value = "'140.64"
$ws.Range('E50').Value = '  +2.07%  '

$ws.Range('D51').Value = "'0.08393"
$ws.Range('E51').Value = '  +5.91%  '
